$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# Update the "time_taken" (F column) timestamps for rows 2-116 in the "data" sheet

$dataWs.Cells.Item(2, 6).Value = "2021-10-05 14:21:28.969004"
$dataWs.Cells.Item(3, 6).Value = "2021-10-05 14:21:28.969012"
$dataWs.Cells.Item(4, 6).Value = "2021-10-05 14:21:28.969016"
$dataWs.Cells.Item(5, 6).Value = "2021-10-05 14:21:28.969018"
$dataWs.Cells.Item(6, 6).Value = "2021-10-05 14:21:28.969021"
$dataWs.Cells.Item(7, 6).Value = "2021-10-05 14:21:28.969024"
$dataWs.Cells.Item(8, 6).Value = "2021-10-05 14:21:28.969026"
$dataWs.Cells.Item(9, 6).Value = "2021-10-05 14:21:28.969029"
$dataWs.Cells.Item(10, 6).Value = "2021-10-05 14:21:28.969032"
$dataWs.Cells.Item(11, 6).Value = "2021-10-05 14:21:28.969034"
$dataWs.Cells.Item(12, 6).Value = "2021-10-05 14:21:28.969037"
$dataWs.Cells.Item(13, 6).Value = "2021-10-05 14:21:28.969039"
$dataWs.Cells.Item(14, 6).Value = "2021-10-05 14:21:28.969042"
$dataWs.Cells.Item(15, 6).Value = "2021-10-05 14:21:28.969044"
$dataWs.Cells.Item(16, 6).Value = "2021-10-05 14:21:28.969047"
$dataWs.Cells.Item(17, 6).Value = "2021-10-05 14:21:28.969050"
$dataWs.Cells.Item(18, 6).Value = "2021-10-05 14:21:28.969052"
$dataWs.Cells.Item(19, 6).Value = "2021-10-05 14:21:28.969055"
$dataWs.Cells.Item(20, 6).Value = "2021-10-05 14:21:28.969058"
$dataWs.Cells.Item(21, 6).Value = "2021-10-05 14:21:28.969077"
$dataWs.Cells.Item(22, 6).Value = "2021-10-05 14:21:28.969083"
$dataWs.Cells.Item(23, 6).Value = "2021-10-05 14:21:28.969085"
$dataWs.Cells.Item(24, 6).Value = "2021-10-05 14:21:28.969088"
$dataWs.Cells.Item(25, 6).Value = "2021-10-05 14:21:28.969090"
$dataWs.Cells.Item(26, 6).Value = "2021-10-05 14:21:28.969093"
$dataWs.Cells.Item(27, 6).Value = "2021-10-05 14:21:28.969096"
$dataWs.Cells.Item(28, 6).Value = "2021-10-05 14:21:28.969098"
$dataWs.Cells.Item(29, 6).Value = "2021-10-05 14:21:28.969101"
$dataWs.Cells.Item(30, 6).Value = "2021-10-05 14:21:28.969103"
$dataWs.Cells.Item(31, 6).Value = "2021-10-05 14:21:28.969106"
$dataWs.Cells.Item(32, 6).Value = "2021-10-05 14:21:28.969109"
$dataWs.Cells.Item(33, 6).Value = "2021-10-05 14:21:28.969111"
$dataWs.Cells.Item(34, 6).Value = "2021-10-05 14:21:28.969116"
$dataWs.Cells.Item(35, 6).Value = "2021-10-05 14:21:28.969118"
$dataWs.Cells.Item(36, 6).Value = "2021-10-05 14:21:28.969121"
$dataWs.Cells.Item(37, 6).Value = "2021-10-05 14:21:28.969124"
$dataWs.Cells.Item(38, 6).Value = "2021-10-05 14:21:28.969126"
$dataWs.Cells.Item(39, 6).Value = "2021-10-05 14:21:28.969128"
$dataWs.Cells.Item(40, 6).Value = "2021-10-05 14:21:28.969131"
$dataWs.Cells.Item(41, 6).Value = "2021-10-05 14:21:28.969133"
$dataWs.Cells.Item(42, 6).Value = "2021-10-05 14:21:28.969136"
$dataWs.Cells.Item(43, 6).Value = "2021-10-05 14:21:28.969176"
$dataWs.Cells.Item(44, 6).Value = "2021-10-05 14:21:28.969182"
$dataWs.Cells.Item(45, 6).Value = "2021-10-05 14:21:28.969185"
$dataWs.Cells.Item(46, 6).Value = "2021-10-05 14:21:28.969188"
$dataWs.Cells.Item(47, 6).Value = "2021-10-05 14:21:28.969190"
$dataWs.Cells.Item(48, 6).Value = "2021-10-05 14:21:28.969193"
$dataWs.Cells.Item(49, 6).Value = "2021-10-05 14:21:28.969196"
$dataWs.Cells.Item(50, 6).Value = "2021-10-05 14:21:28.969198"
$dataWs.Cells.Item(51, 6).Value = "2021-10-05 14:21:28.969201"
$dataWs.Cells.Item(52, 6).Value = "2021-10-05 14:21:28.969204"
$dataWs.Cells.Item(53, 6).Value = "2021-10-05 14:21:28.969207"
$dataWs.Cells.Item(54, 6).Value = "2021-10-05 14:21:28.969210"
$dataWs.Cells.Item(55, 6).Value = "2021-10-05 14:21:28.969213"
$dataWs.Cells.Item(56, 6).Value = "2021-10-05 14:21:28.969215"
$dataWs.Cells.Item(57, 6).Value = "2021-10-05 14:21:28.969218"
$dataWs.Cells.Item(58, 6).Value = "2021-10-05 14:21:28.969220"
$dataWs.Cells.Item(59, 6).Value = "2021-10-05 14:21:28.969223"
$dataWs.Cells.Item(60, 6).Value = "2021-10-05 14:21:28.969226"
$dataWs.Cells.Item(61, 6).Value = "2021-10-05 14:21:28.969228"
$dataWs.Cells.Item(62, 6).Value = "2021-10-05 14:21:28.969231"
$dataWs.Cells.Item(63, 6).Value = "2021-10-05 14:21:28.969233"
$dataWs.Cells.Item(64, 6).Value = "2021-10-05 14:21:28.969236"
$dataWs.Cells.Item(65, 6).Value = "2021-10-05 14:21:28.969238"
$dataWs.Cells.Item(66, 6).Value = "2021-10-05 14:21:28.969242"
$dataWs.Cells.Item(67, 6).Value = "2021-10-05 14:21:28.969245"
$dataWs.Cells.Item(68, 6).Value = "2021-10-05 14:21:28.969248"
$dataWs.Cells.Item(69, 6).Value = "2021-10-05 14:21:28.969251"
$dataWs.Cells.Item(70, 6).Value = "2021-10-05 14:21:28.969253"
$dataWs.Cells.Item(71, 6).Value = "2021-10-05 14:21:28.969256"
$dataWs.Cells.Item(72, 6).Value = "2021-10-05 14:21:28.969259"
$dataWs.Cells.Item(73, 6).Value = "2021-10-05 14:21:28.969261"
$dataWs.Cells.Item(74, 6).Value = "2021-10-05 14:21:28.969264"
$dataWs.Cells.Item(75, 6).Value = "2021-10-05 14:21:28.969266"
$dataWs.Cells.Item(76, 6).Value = "2021-10-05 14:21:28.969269"
$dataWs.Cells.Item(77, 6).Value = "2021-10-05 14:21:28.969272"
$dataWs.Cells.Item(78, 6).Value = "2021-10-05 14:21:28.969276"
$dataWs.Cells.Item(79, 6).Value = "2021-10-05 14:21:28.969279"
$dataWs.Cells.Item(80, 6).Value = "2021-10-05 14:21:28.969282"
$dataWs.Cells.Item(81, 6).Value = "2021-10-05 14:21:28.969285"
$dataWs.Cells.Item(82, 6).Value = "2021-10-05 14:21:28.969288"
$dataWs.Cells.Item(83, 6).Value = "2021-10-05 14:21:28.969290"
$dataWs.Cells.Item(84, 6).Value = "2021-10-05 14:21:28.969293"
$dataWs.Cells.Item(85, 6).Value = "2021-10-05 14:21:28.969295"
$dataWs.Cells.Item(86, 6).Value = "2021-10-05 14:21:28.969298"
$dataWs.Cells.Item(87, 6).Value = "2021-10-05 14:21:28.969301"
$dataWs.Cells.Item(88, 6).Value = "2021-10-05 14:21:28.969303"
$dataWs.Cells.Item(89, 6).Value = "2021-10-05 14:21:28.969306"
$dataWs.Cells.Item(90, 6).Value = "2021-10-05 14:21:28.969309"
$dataWs.Cells.Item(91, 6).Value = "2021-10-05 14:21:28.969312"
$dataWs.Cells.Item(92, 6).Value = "2021-10-05 14:21:28.969314"
$dataWs.Cells.Item(93, 6).Value = "2021-10-05 14:21:28.969317"
$dataWs.Cells.Item(94, 6).Value = "2021-10-05 14:21:28.969321"
$dataWs.Cells.Item(95, 6).Value = "2021-10-05 14:21:28.969323"
$dataWs.Cells.Item(96, 6).Value = "2021-10-05 14:21:28.969326"
$dataWs.Cells.Item(97, 6).Value = "2021-10-05 14:21:28.969329"
$dataWs.Cells.Item(98, 6).Value = "2021-10-05 14:21:28.969332"
$dataWs.Cells.Item(99, 6).Value = "2021-10-05 14:21:28.969334"
$dataWs.Cells.Item(100, 6).Value = "2021-10-05 14:21:28.969337"
$dataWs.Cells.Item(101, 6).Value = "2021-10-05 14:21:28.969339"
$dataWs.Cells.Item(102, 6).Value = "2021-10-05 14:21:28.969342"
$dataWs.Cells.Item(103, 6).Value = "2021-10-05 14:21:28.969345"
$dataWs.Cells.Item(104, 6).Value = "2021-10-05 14:21:28.969348"
$dataWs.Cells.Item(105, 6).Value = "2021-10-05 14:21:28.969350"
$dataWs.Cells.Item(106, 6).Value = "2021-10-05 14:21:28.969353"
$dataWs.Cells.Item(107, 6).Value = "2021-10-05 14:21:28.969356"
$dataWs.Cells.Item(108, 6).Value = "2021-10-05 14:21:28.969358"
$dataWs.Cells.Item(109, 6).Value = "2021-10-05 14:21:28.969361"
$dataWs.Cells.Item(110, 6).Value = "2021-10-05 14:21:28.969366"
$dataWs.Cells.Item(111, 6).Value = "2021-10-05 14:21:28.969369"
$dataWs.Cells.Item(112, 6).Value = "2021-10-05 14:21:28.969371"
$dataWs.Cells.Item(113, 6).Value = "2021-10-05 14:21:28.969374"
$dataWs.Cells.Item(114, 6).Value = "2021-10-05 14:21:28.969376"
$dataWs.Cells.Item(115, 6).Value = "2021-10-05 14:21:28.969379"
$dataWs.Cells.Item(116, 6).Value = "2021-10-05 14:21:28.969382"

# Add a new "metadata" sheet positioned right after "data"
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Reuse the header formatting (bold, bordered, centered) from the "data" sheet's
# header row, and the index-column formatting from its first data row, so the
# new sheet matches the workbook's existing look.
$dataWs.Range("B1:F1").Copy()
$metaWs.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dataWs.Range("A2").Copy()
$metaWs.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$metaWs.Cells.Item(1, 2).Value = "data_name"
$metaWs.Cells.Item(1, 3).Value = "data_id"
$metaWs.Cells.Item(1, 4).Value = "data_version"
$metaWs.Cells.Item(1, 5).Value = "data_version_created"
$metaWs.Cells.Item(1, 6).Value = "panel_query_time"
$metaWs.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$metaWs.Cells.Item(2, 1).Value = 0
$metaWs.Cells.Item(2, 2).Value = "Malformations of cortical development"
$metaWs.Cells.Item(2, 3).Value = 96

# data_version "2.92" must stay text, not become the number 2.92
$metaWs.Cells.Item(2, 4).NumberFormat = "@"
$metaWs.Cells.Item(2, 4).Value = "2.92"
$metaWs.Cells.Item(2, 4).Style = "Normal"

$metaWs.Cells.Item(2, 5).Value = "2021-10-01T14:43:12.917462Z"
$metaWs.Cells.Item(2, 6).Value = "2021-10-05 14:21:28.965814"
$metaWs.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/96/?format=json"

